$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.825.69"
$ws.Range("E2").Value = "  +3.39%  "

$ws.Range("D3").Value = "3.133.56"
$ws.Range("E3").Value = "  +2.40%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.09"
$ws.Range("E5").Value = "  +2.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.69"
$ws.Range("E6").Value = "  +2.84%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.124.05"
$ws.Range("E8").Value = "  +2.46%  "

$ws.Range("E9").Value = "  +2.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +19.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.70"
$ws.Range("E11").Value = "  +3.99%  "

$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("E13").Value = "  +7.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.03"
$ws.Range("E14").Value = "  +3.87%  "

$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("D16").Value = "3.647.67"
$ws.Range("E16").Value = "  +2.32%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.17"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "63.717.69"
$ws.Range("E18").Value = "  +3.27%  "

$ws.Range("D19").Value = "3.128.78"
$ws.Range("E19").Value = "  +2.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.00"
$ws.Range("E20").Value = "  +3.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.22"
$ws.Range("E21").Value = "  +2.31%  "

$ws.Range("E22").Value = "  +0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.54"
$ws.Range("E23").Value = "  +3.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.29"
$ws.Range("E24").Value = "  -2.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.34"
$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.66"
$ws.Range("E27").Value = "  +8.48%  "

$ws.Range("E28").Value = "  +3.06%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.86"
$ws.Range("E31").Value = "  +4.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.04"
$ws.Range("E32").Value = "  +2.03%  "

$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("D34").Value = "0.0₃0871"
$ws.Range("E34").Value = "  +7.90%  "

$ws.Range("E35").Value = "  +10.16%  "

$ws.Range("E36").Value = "  +1.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("E37").Value = "  +15.76%  "

$ws.Range("E38").Value = "  +1.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.81"
$ws.Range("E39").Value = "  +1.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "451.21"
$ws.Range("E40").Value = "  +9.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.71"
$ws.Range("E41").Value = "  -1.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0374"
$ws.Range("E42").Value = "  +1.97%  "

$ws.Range("D43").Value = "2.919.70"
$ws.Range("E43").Value = "  +5.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.277"
$ws.Range("E44").Value = "  +5.36%  "

$ws.Range("E45").Value = "  +3.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  +3.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.10"
$ws.Range("E47").Value = "  +4.16%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.45"
$ws.Range("E49").Value = "  -6.71%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.111"
$ws.Range("E50").Value = "  +0.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.71"
$ws.Range("E51").Value = "  +2.89%  "
